$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Q7)
$ws.Range("B9").Value = 0.06204647633618569
$ws.Range("C9").Value = 0.2035131537503989
$ws.Range("D9").Value = 0.04871313030636959
$ws.Range("E9").Value = 0.2207105124509696
$ws.Range("F9").Value = 0.2320259427235655
$ws.Range("G9").Value = 6

# Row 10 (Q8)
$ws.Range("B10").Value = 0.2179584722456556
$ws.Range("C10").Value = 0.2179584722456556
$ws.Range("D10").Value = 0.1208316223875407
$ws.Range("E10").Value = 0.3476084325610365
$ws.Range("F10").Value = 0.3316452775870941
$ws.Range("G10").Value = 3

# Row 11 (Q9)
$ws.Range("B11").Value = 0.3403795785247692
$ws.Range("C11").Value = 0.3403795785247692
$ws.Range("D11").Value = 0.1158582574766995
$ws.Range("E11").Value = 0.3403795785247692
$ws.Range("F11").ClearContents()
$ws.Range("G11").Value = 1
